$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29 (pushes existing rows 29..78 down to 30..79)
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly data point
$ws.Range("A29").Value = 2
$ws.Range("B29").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 45070
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 100112032
$ws.Range("G29").Value = "Zapallo italiano"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 700
$ws.Range("K29").Value = 7000
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 7500
$ws.Range("N29").Value = "$/caja 60 unidades"
$ws.Range("O29").Value = "Provincia de Limarí"
$ws.Range("P29").Value = 125
$ws.Range("Q29").Value = 60
$ws.Range("R29").Value = "Hortaliza"
